# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stock) sheet (Worksheets index 5) gets three new trailing
# columns: date / legislator_name / legislator_id, mirroring the header
# style of the existing columns (bold, centered, thin border) on row 1
# and the plain data style on row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# ---- Header row (row 1): H1=date, I1=legislator_name, J1=legislator_id
# Apply the same look as the existing header cells (bold, centered,
# thin box border) BEFORE writing the value - cell by cell, so each one
# resolves to the workbook's existing shared header style.
foreach ($col in @("H", "I", "J")) {
    $hdr = $ws.Range($col + "1")
    $hdr.Font.Bold = $true
    $hdr.HorizontalAlignment = -4108   # xlCenter
    $hdr.VerticalAlignment = -4160     # xlTop
    $hdr.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $hdr.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $hdr.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $hdr.Borders.Item(10).LineStyle = 1  # xlEdgeRight
}
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# ---- Data row (row 2): H2=date value, I2=legislator name, J2=legislator id
# Write the date as literal text (not an auto-converted date serial) by
# computing it through TEXT() and pasting the result back as a value.
$dateCell = $ws.Range("H2")
$dateCell.Formula = '=TEXT("2012-03-26","@")'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)  # xlPasteValues

$ws.Range("I2").Value = "潘維剛"
$ws.Range("J2").Value = 678
